$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.724724292755127
$ws.Range("B1").Value = 1.458641648292542
$ws.Range("D1").Value = 1.96343207359314
$ws.Range("E1").Value = 1.232074737548828
